$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "71.770.67"
Set-TextValue "E2" "  -1.63%  "

# Row 3
Set-TextValue "D3" "3.886.91"
Set-TextValue "E3" "  -2.59%  "

# Row 4
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.34%  "

# Row 5
Set-TextValue "D5" "600.17"
Set-TextValue "E5" "  +1.64%  "

# Row 6
Set-TextValue "D6" "169.20"
Set-TextValue "E6" "  +7.85%  "

# Row 7
Set-TextValue "D7" "0.677"
Set-TextValue "E7" "  -0.42%  "

# Row 8
Set-TextValue "E8" "  +0.19%  "

# Row 9
Set-TextValue "D9" "0.759"
Set-TextValue "E9" "  +1.56%  "

# Row 10
Set-TextValue "D10" "0.178"
Set-TextValue "E10" "  +5.92%  "

# Row 11
Set-TextValue "D11" "53.78"
Set-TextValue "E11" "  +0.73%  "

# Row 12
Set-TextValue "D12" "0.0000322"
Set-TextValue "E12" "  +1.35%  "

# Row 13
Set-TextValue "D13" "11.28"
Set-TextValue "E13" "  +4.25%  "

# Row 14
Set-TextValue "D14" "4.525.00"
Set-TextValue "E14" "  -1.81%  "

# Row 15
Set-TextValue "D15" "3.883.33"
Set-TextValue "E15" "  -1.82%  "

# Row 16
Set-TextValue "D16" "20.95"
Set-TextValue "E16" "  +2.38%  "

# Row 17
Set-TextValue "D17" "13.93"
Set-TextValue "E17" "  -0.85%  "

# Row 18
Set-TextValue "E18" "  -4.85%  "

# Row 19
Set-TextValue "E19" "  -1.92%  "

# Row 20
Set-TextValue "D20" "71.667.24"
Set-TextValue "E20" "  -1.03%  "

# Row 21
Set-TextValue "D21" "436.12"
Set-TextValue "E21" "  +1.41%  "

# Row 22
Set-TextValue "D22" "4.77"
Set-TextValue "E22" "  +1.96%  "

# Row 23
Set-TextValue "D23" "94.78"
Set-TextValue "E23" "  -1.19%  "

# Row 24
Set-TextValue "D24" "3.30"
Set-TextValue "E24" "  -3.46%  "

# Row 25
Set-TextValue "D25" "13.87"
Set-TextValue "E25" "  -3.41%  "

# Row 26
Set-TextValue "D26" "4.11"
Set-TextValue "E26" "  -6.31%  "

# Row 27
Set-TextValue "D27" "10.99"
Set-TextValue "E27" "  -3.31%  "

# Row 28
Set-TextValue "D28" "5.95"
Set-TextValue "E28" "  +0.50%  "

# Row 29
Set-TextValue "D29" "10.22"
Set-TextValue "E29" "  -5.72%  "

# Row 30
Set-TextValue "D30" "35.18"
Set-TextValue "E30" "  -3.39%  "

# Row 31
Set-TextValue "D31" "7.87"
Set-TextValue "E31" "  +0.47%  "

# Row 32
Set-TextValue "D32" "51.34"
Set-TextValue "E32" "  +1.51%  "

# Row 33
Set-TextValue "D33" "13.61"
Set-TextValue "E33" "  +0.87%  "

# Row 34
Set-TextValue "E34" "  -4.08%  "

# Row 35
Set-TextValue "D35" "0.0₃0988"
Set-TextValue "E35" "  +15.20%  "

# Row 36
Set-TextValue "D36" "69.01"
Set-TextValue "E36" "  -0.77%  "

# Row 37
Set-TextValue "D37" "618.55"
Set-TextValue "E37" "  -9.38%  "

# Row 38
Set-TextValue "D38" "0.421"
Set-TextValue "E38" "  -4.22%  "

# Row 39
Set-TextValue "D39" "1.00"
Set-TextValue "E39" "  -0.05%  "

# Row 40
Set-TextValue "D40" "3.30"
Set-TextValue "E40" "  -0.54%  "

# Row 41
Set-TextValue "E41" "  -1.26%  "

# Row 42
Set-TextValue "D42" "0.999"
Set-TextValue "E42" "  +0.10%  "

# Row 43
Set-TextValue "D43" "3.22"
Set-TextValue "E43" "  +34.17%  "

# Row 44
Set-TextValue "D44" "0.0471"
Set-TextValue "E44" "  -2.96%  "

# Row 45
Set-TextValue "E45" "  -6.37%  "

# Row 46
Set-TextValue "D46" "2.65"
Set-TextValue "E46" "  -2.38%  "

# Row 47
Set-TextValue "D47" "0.144"
Set-TextValue "E47" "  -2.72%  "

# Row 48
Set-TextValue "E48" "  -0.69%  "

# Row 51
Set-TextValue "D51" "0.000272"
Set-TextValue "E51" "  +0.78%  "

# Row 49 and 50 swap (Maker <-> WEMIXToken)
Set-TextValue "B49" "WEMIXToken"
Set-TextValue "C49" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D49" "2.78"
Set-TextValue "E49" "  -17.43%  "

Set-TextValue "B50" "Maker"
Set-TextValue "C50" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D50" "2.846.14"
Set-TextValue "E50" "  +2.66%  "

